$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "271.32"
Set-TextValue "G2" "15"
Set-TextValue "D3" "22.98"
Set-TextValue "G3" "15"
Set-TextValue "D4" "6.386"
Set-TextValue "G4" "15"
Set-TextValue "D5" "0.06279"
Set-TextValue "G5" "15"
Set-TextValue "D6" "3.652"
Set-TextValue "G6" "15"
Set-TextValue "D7" "6.713"
Set-TextValue "G7" "15"
Set-TextValue "D8" "1.385"
Set-TextValue "G8" "15"
Set-TextValue "D9" "0.8366"
Set-TextValue "G9" "15"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1635"
$ws.Range("E10").Value = "9WazirXWRX"
Set-TextValue "G10" "15"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.08408"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
Set-TextValue "G11" "15"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D12" "0.03489"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue "G12" "15"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.03141"
$ws.Range("E13").Value = "12BitrueCoinBTR"
Set-TextValue "G13" "15"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D14" "0.09336"
$ws.Range("E14").Value = "13BitMartTokenBMX"
Set-TextValue "G14" "15"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D15" "3.892"
$ws.Range("E15").Value = "14MCDexMCB"
Set-TextValue "G15" "15"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D16" "0.001700"
$ws.Range("E16").Value = "15BitForexTokenBF"
Set-TextValue "G16" "15"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04850"
$ws.Range("E17").Value = "16CoinExTokenCET"
Set-TextValue "G17" "15"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D18" "0.006273"
$ws.Range("E18").Value = "17TigerCashTCH"
Set-TextValue "G18" "15"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D19" "0.004681"
$ws.Range("E19").Value = "18HotbitTokenHTBWorstin24h"
Set-TextValue "G19" "15"
Set-TextValue "G20" "15"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D21" "0.0001497"
$ws.Range("E21").Value = "20NitroExNTX"
Set-TextValue "G21" "15"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D22" "3.741"
$ws.Range("E22").Value = "21LEOLEO"
Set-TextValue "G22" "15"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D23" "2.323"
$ws.Range("E23").Value = "22BTSETokenBTSE"
Set-TextValue "G23" "15"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D24" "0.01380"
$ws.Range("E24").Value = "23OneONE"
Set-TextValue "G24" "15"
Set-TextValue "D25" "0.3406"
Set-TextValue "G25" "15"
Set-TextValue "D26" "0.1261"
Set-TextValue "G26" "15"
Set-TextValue "G27" "15"
Set-TextValue "G28" "15"
Set-TextValue "G29" "15"
Set-TextValue "G30" "15"
Set-TextValue "G31" "15"
Set-TextValue "G32" "15"
Set-TextValue "G33" "15"
Set-TextValue "G34" "15"
Set-TextValue "G35" "15"
Set-TextValue "G36" "15"
Set-TextValue "G37" "15"
Set-TextValue "G38" "15"
Set-TextValue "G39" "15"
Set-TextValue "D40" "0.04694"
Set-TextValue "G40" "15"
Set-TextValue "D41" "0.006887"
Set-TextValue "G41" "15"
Set-TextValue "D42" "0.1178"
Set-TextValue "G42" "15"
Set-TextValue "D43" "0.003449"
Set-TextValue "G43" "15"
Set-TextValue "G44" "15"
Set-TextValue "D45" "0.00006246"
Set-TextValue "G45" "15"
Set-TextValue "D46" "0.00000000749"
Set-TextValue "G46" "15"
Set-TextValue "D47" "0.7970"
Set-TextValue "G47" "15"
Set-TextValue "D48" "0.09212"
Set-TextValue "G48" "15"
Set-TextValue "G49" "15"
Set-TextValue "D50" "0.01238"
Set-TextValue "G50" "15"
Set-TextValue "G51" "15"
